$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" (D) and "Volume(1h)" (E) columns hold text values such as
# "64.348.28" (thousands separated by dots), "163.00" or padded percentage
# strings like "  -0.45%  ". Force a text number format on every cell we
# touch so Excel keeps the exact string instead of coercing it to a number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D2").Value2 = '64.348.28'
$ws.Range("E2").Value2 = '  -0.45%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D3").Value2 = '3.132.96'
$ws.Range("E3").Value2 = '  -1.71%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value2 = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D5").Value2 = '571.99'
$ws.Range("E5").Value2 = '  -0.20%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D6").Value2 = '163.93'
$ws.Range("E6").Value2 = '  -3.95%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D7").Value2 = '0.999'
$ws.Range("E7").Value2 = '  -0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D8").Value2 = '0.573'
$ws.Range("E8").Value2 = '  -6.19%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D9").Value2 = '3.146.20'
$ws.Range("E9").Value2 = '  -1.58%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value2 = '  -3.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D11").Value2 = '6.63'
$ws.Range("E11").Value2 = '  -3.49%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D12").Value2 = '0.383'
$ws.Range("E12").Value2 = '  -2.52%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D13").Value2 = '3.685.88'
$ws.Range("E13").Value2 = '  -1.73%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D14").Value2 = '0.126'
$ws.Range("E14").Value2 = '  -2.11%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D15").Value2 = '64.406.02'
$ws.Range("E15").Value2 = '  -0.40%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D16").Value2 = '24.87'
$ws.Range("E16").Value2 = '  -2.75%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D17").Value2 = '3.149.33'
$ws.Range("E17").Value2 = '  -1.30%  '

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value2 = '  -2.76%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D19").Value2 = '406.41'
$ws.Range("E19").Value2 = '  -3.76%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D20").Value2 = '5.23'
$ws.Range("E20").Value2 = '  -2.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D21").Value2 = '12.49'
$ws.Range("E21").Value2 = '  -4.19%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D22").Value2 = '7.04'
$ws.Range("E22").Value2 = '  -1.91%  '

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value2 = '  -0.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D24").Value2 = '68.55'
$ws.Range("E24").Value2 = '  -2.70%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D25").Value2 = '0.482'
$ws.Range("E25").Value2 = '  -3.93%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("D26").Value2 = '0.193'
$ws.Range("E26").Value2 = '  -6.32%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D27").Value2 = '0.0000101'
$ws.Range("E27").Value2 = '  -4.60%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D28").Value2 = '8.86'
$ws.Range("E28").Value2 = '  +0.19%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("D29").Value2 = '0.996'
$ws.Range("E29").Value2 = '  +0.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D30").Value2 = '0.999'
$ws.Range("E30").Value2 = '  +0.07%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("D31").Value2 = '1.80'
$ws.Range("E31").Value2 = '  -2.06%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D32").Value2 = '21.22'
$ws.Range("E32").Value2 = '  -3.06%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("D33").Value2 = '163.00'
$ws.Range("E33").Value2 = '  +3.80%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("D34").Value2 = '4.83'
$ws.Range("E34").Value2 = '  -5.07%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("D35").Value2 = '6.26'
$ws.Range("E35").Value2 = '  -2.34%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("D36").Value2 = '1.12'
$ws.Range("E36").Value2 = '  -0.68%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D37").Value2 = '1.35'
$ws.Range("E37").Value2 = '  -1.35%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value2 = '  -2.23%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D39").Value2 = '2.633.28'
$ws.Range("E39").Value2 = '  -3.24%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D40").Value2 = '23.62'
$ws.Range("E40").Value2 = '  -3.54%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D41").Value2 = '4.09'
$ws.Range("E41").Value2 = '  -4.31%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D42").Value2 = '38.19'
$ws.Range("E42").Value2 = '  -2.57%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D43").Value2 = '0.691'
$ws.Range("E43").Value2 = '  -4.22%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D44").Value2 = '0.0615'
$ws.Range("E44").Value2 = '  -1.73%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value2 = '  -5.41%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("B46").Value2 = 'InjectiveProtocol'
$ws.Range("C46").Value2 = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").Value2 = '21.32'
$ws.Range("E46").Value2 = '  -1.57%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("B47").Value2 = 'Bittensor'
$ws.Range("C47").Value2 = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").Value2 = '288.97'
$ws.Range("E47").Value2 = '  -1.73%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D48").Value2 = '0.0253'
$ws.Range("E48").Value2 = '  -3.96%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D49").Value2 = '0.997'
$ws.Range("E49").Value2 = '  -0.05%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value2 = '  -1.90%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value2 = '  +0.57%  '
